# The workbook tracks daily market-price records for "Pepino ensalada" at
# Feria Lagunitas de Puerto Montt. A new weekly record is inserted as a new
# row 37, shifting all existing records (previously rows 37-156) down by one
# row (to rows 38-157). The sheet's used range grows from A1:R156 to A1:R157.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 37, shifting rows 37:156 down to 38:157.
$ws.Rows("37:37").Insert()

# Populate the newly inserted row 37 with the new record's data.
$ws.Range("A37").Value = 4
$ws.Range("B37").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C37").Value = "Los Lagos"
$ws.Range("D37").Value = 44481
$ws.Range("E37").Value = 10
$ws.Range("F37").Value = 100112043
$ws.Range("G37").Value = "Pepino ensalada"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 400
$ws.Range("K37").Value = 16000
$ws.Range("L37").Value = 17000
$ws.Range("M37").Value = 16500
$ws.Range("N37").Value = "$/caja 60 unidades"
$ws.Range("O37").Value = "Región de Arica y Parinacota"
$ws.Range("P37").Value = 275
$ws.Range("Q37").Value = 60
$ws.Range("R37").Value = "Hortaliza"
